$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - extraction-template-12.xlsx
$ws.Range("B2").Value = "\ExtractionTemplate\ManagePopulations\extraction-template-12.xlsx"
$ws.Range("A2").Value = "extraction-template-12.xlsx"
$ws.Range("A2:B2").Style = "Normal"

# Row 3 - extraction-template-14.xlsx
$ws.Range("B3").Value = "\ExtractionTemplate\ManagePopulations\extraction-template-14.xlsx"
$ws.Range("A3").Value = "extraction-template-14.xlsx"

# Row 4 - extraction-template-17.xlsx
$ws.Range("B4").Value = "\ExtractionTemplate\ManagePopulations\extraction-template-17.xlsx"
$ws.Range("A4").Value = "extraction-template-17.xlsx"

# Recompute best-fit column widths based on new content (matches original AutoFit)
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Update selection to B4
$ws.Range("B4").Select()
